$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 374-393 with corrected weekly price data
# Row 374
$ws.Cells.Item(374, 4).Value = 44585
$ws.Cells.Item(374, 10).Value = 1200
$ws.Cells.Item(374, 15).Value = 'Región de Los Lagos'

# Row 375
$ws.Cells.Item(375, 4).Value = 44585
$ws.Cells.Item(375, 9).Value = '1a nueva(o)'
$ws.Cells.Item(375, 11).Value = 7000
$ws.Cells.Item(375, 12).Value = 7000
$ws.Cells.Item(375, 13).Value = 7000
$ws.Cells.Item(375, 16).Value = 280

# Row 376
$ws.Cells.Item(376, 4).Value = 44560
$ws.Cells.Item(376, 8).Value = 'Asterix'
$ws.Cells.Item(376, 9).Value = '1a nueva(o)'
$ws.Cells.Item(376, 10).Value = 1600
$ws.Cells.Item(376, 11).Value = 7000
$ws.Cells.Item(376, 12).Value = 7000
$ws.Cells.Item(376, 13).Value = 7000
$ws.Cells.Item(376, 15).Value = 'Región del Maule'
$ws.Cells.Item(376, 16).Value = 280

# Row 377
$ws.Cells.Item(377, 4).Value = 44272
$ws.Cells.Item(377, 8).Value = 'Patagonia'
$ws.Cells.Item(377, 9).Value = '1a (cosecha)'
$ws.Cells.Item(377, 10).Value = 1200
$ws.Cells.Item(377, 11).Value = 6000
$ws.Cells.Item(377, 12).Value = 6000
$ws.Cells.Item(377, 13).Value = 6000
$ws.Cells.Item(377, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(377, 16).Value = 240

# Row 378
$ws.Cells.Item(378, 4).Value = 44272
$ws.Cells.Item(378, 8).Value = 'Rodeo'
$ws.Cells.Item(378, 9).Value = '1a (cosecha)'
$ws.Cells.Item(378, 10).Value = 1300
$ws.Cells.Item(378, 11).Value = 6000
$ws.Cells.Item(378, 12).Value = 6000
$ws.Cells.Item(378, 13).Value = 6000
$ws.Cells.Item(378, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(378, 16).Value = 240

# Row 379
$ws.Cells.Item(379, 4).Value = 44162
$ws.Cells.Item(379, 8).Value = 'Rosara'
$ws.Cells.Item(379, 10).Value = 2400
$ws.Cells.Item(379, 12).Value = 9000
$ws.Cells.Item(379, 13).Value = 8500
$ws.Cells.Item(379, 16).Value = 340

# Row 380
$ws.Cells.Item(380, 4).Value = 44529
$ws.Cells.Item(380, 9).Value = '1a nueva(o)'
$ws.Cells.Item(380, 11).Value = 9000
$ws.Cells.Item(380, 12).Value = 9000
$ws.Cells.Item(380, 13).Value = 9000
$ws.Cells.Item(380, 16).Value = 360

# Row 381
$ws.Cells.Item(381, 4).Value = 44529
$ws.Cells.Item(381, 9).Value = '1a nueva(o)'
$ws.Cells.Item(381, 10).Value = 1200
$ws.Cells.Item(381, 11).Value = 8000
$ws.Cells.Item(381, 12).Value = 8000
$ws.Cells.Item(381, 13).Value = 8000
$ws.Cells.Item(381, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(381, 15).Value = 'Región del Maule'
$ws.Cells.Item(381, 16).Value = 320

# Row 382
$ws.Cells.Item(382, 4).Value = 44414
$ws.Cells.Item(382, 8).Value = 'Asterix'
$ws.Cells.Item(382, 9).Value = '1a (guarda)'
$ws.Cells.Item(382, 11).Value = 6000
$ws.Cells.Item(382, 12).Value = 6000
$ws.Cells.Item(382, 13).Value = 6000
$ws.Cells.Item(382, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(382, 15).Value = 'Región del Maule'
$ws.Cells.Item(382, 16).Value = 240

# Row 383
$ws.Cells.Item(383, 4).Value = 44414
$ws.Cells.Item(383, 8).Value = 'Rodeo'
$ws.Cells.Item(383, 9).Value = '1a (guarda lavada)'
$ws.Cells.Item(383, 10).Value = 1500
$ws.Cells.Item(383, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(383, 15).Value = 'Provincia de Valdivia'

# Row 384
$ws.Cells.Item(384, 4).Value = 44323
$ws.Cells.Item(384, 9).Value = '1a (cosecha lavada)'
$ws.Cells.Item(384, 11).Value = 6500
$ws.Cells.Item(384, 12).Value = 6500
$ws.Cells.Item(384, 13).Value = 6500
$ws.Cells.Item(384, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(384, 16).Value = 260

# Row 385
$ws.Cells.Item(385, 4).Value = 44306
$ws.Cells.Item(385, 8).Value = 'Patagonia'
$ws.Cells.Item(385, 9).Value = '1a (cosecha)'
$ws.Cells.Item(385, 10).Value = 500
$ws.Cells.Item(385, 11).Value = 6000
$ws.Cells.Item(385, 12).Value = 6000
$ws.Cells.Item(385, 13).Value = 6000
$ws.Cells.Item(385, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(385, 15).Value = 'Región del Maule'
$ws.Cells.Item(385, 16).Value = 240

# Row 386
$ws.Cells.Item(386, 4).Value = 44299
$ws.Cells.Item(386, 11).Value = 6000
$ws.Cells.Item(386, 12).Value = 6000
$ws.Cells.Item(386, 13).Value = 6000
$ws.Cells.Item(386, 16).Value = 240

# Row 387
$ws.Cells.Item(387, 4).Value = 44428
$ws.Cells.Item(387, 8).Value = 'Rodeo'
$ws.Cells.Item(387, 9).Value = '1a (guarda lavada)'
$ws.Cells.Item(387, 10).Value = 1500
$ws.Cells.Item(387, 11).Value = 7000
$ws.Cells.Item(387, 12).Value = 7000
$ws.Cells.Item(387, 13).Value = 7000
$ws.Cells.Item(387, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(387, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(387, 16).Value = 280

# Row 388
$ws.Cells.Item(388, 4).Value = 44321
$ws.Cells.Item(388, 8).Value = 'Patagonia'
$ws.Cells.Item(388, 9).Value = '1a (cosecha)'
$ws.Cells.Item(388, 10).Value = 1200
$ws.Cells.Item(388, 11).Value = 5500
$ws.Cells.Item(388, 12).Value = 5500
$ws.Cells.Item(388, 13).Value = 5500
$ws.Cells.Item(388, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(388, 16).Value = 220

# Row 389
$ws.Cells.Item(389, 4).Value = 44302
$ws.Cells.Item(389, 8).Value = 'Asterix'
$ws.Cells.Item(389, 9).Value = '1a (cosecha)'
$ws.Cells.Item(389, 10).Value = 1200
$ws.Cells.Item(389, 11).Value = 5500
$ws.Cells.Item(389, 12).Value = 5500
$ws.Cells.Item(389, 13).Value = 5500
$ws.Cells.Item(389, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(389, 16).Value = 220

# Row 390
$ws.Cells.Item(390, 4).Value = 44209
$ws.Cells.Item(390, 8).Value = 'Asterix'
$ws.Cells.Item(390, 9).Value = '1a nueva(o)'
$ws.Cells.Item(390, 10).Value = 800
$ws.Cells.Item(390, 11).Value = 10000
$ws.Cells.Item(390, 12).Value = 10000
$ws.Cells.Item(390, 13).Value = 10000
$ws.Cells.Item(390, 15).Value = 'Región del Maule'
$ws.Cells.Item(390, 16).Value = 400

# Row 391
$ws.Cells.Item(391, 4).Value = 44209
$ws.Cells.Item(391, 8).Value = 'Rosara'
$ws.Cells.Item(391, 9).Value = '1a nueva(o)'
$ws.Cells.Item(391, 10).Value = 800
$ws.Cells.Item(391, 11).Value = 8000
$ws.Cells.Item(391, 12).Value = 8000
$ws.Cells.Item(391, 13).Value = 8000
$ws.Cells.Item(391, 15).Value = 'Región del Maule'
$ws.Cells.Item(391, 16).Value = 320

# Row 392
$ws.Cells.Item(392, 4).Value = 44274
$ws.Cells.Item(392, 8).Value = 'Patagonia'
$ws.Cells.Item(392, 9).Value = '1a (cosecha)'
$ws.Cells.Item(392, 11).Value = 6000
$ws.Cells.Item(392, 12).Value = 6000
$ws.Cells.Item(392, 13).Value = 6000
$ws.Cells.Item(392, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(392, 16).Value = 240

# Row 393
$ws.Cells.Item(393, 4).Value = 44274
$ws.Cells.Item(393, 9).Value = '1a (cosecha)'
$ws.Cells.Item(393, 11).Value = 6000
$ws.Cells.Item(393, 12).Value = 6000
$ws.Cells.Item(393, 13).Value = 6000
$ws.Cells.Item(393, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(393, 16).Value = 240

# Row 394
$ws.Cells.Item(394, 1).Value = 5
$ws.Cells.Item(394, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(394, 3).Value = 'Maule'
$ws.Cells.Item(394, 4).Value = 44554
$ws.Cells.Item(394, 5).Value = 7
$ws.Cells.Item(394, 6).Value = 100114001
$ws.Cells.Item(394, 7).Value = 'Papa'
$ws.Cells.Item(394, 8).Value = 'Asterix'
$ws.Cells.Item(394, 9).Value = '1a nueva(o)'
$ws.Cells.Item(394, 10).Value = 1200
$ws.Cells.Item(394, 11).Value = 9000
$ws.Cells.Item(394, 12).Value = 9000
$ws.Cells.Item(394, 13).Value = 9000
$ws.Cells.Item(394, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(394, 15).Value = 'Región del Maule'
$ws.Cells.Item(394, 16).Value = 360
$ws.Cells.Item(394, 17).Value = 25
$ws.Cells.Item(394, 18).Value = 'Hortaliza'
$ws.Cells.Item(394, 4).NumberFormat = $ws.Cells.Item(373, 4).NumberFormat

# Row 395
$ws.Cells.Item(395, 1).Value = 5
$ws.Cells.Item(395, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(395, 3).Value = 'Maule'
$ws.Cells.Item(395, 4).Value = 44392
$ws.Cells.Item(395, 5).Value = 7
$ws.Cells.Item(395, 6).Value = 100114001
$ws.Cells.Item(395, 7).Value = 'Papa'
$ws.Cells.Item(395, 8).Value = 'Rodeo'
$ws.Cells.Item(395, 9).Value = '1a (cosecha lavada)'
$ws.Cells.Item(395, 10).Value = 1200
$ws.Cells.Item(395, 11).Value = 7000
$ws.Cells.Item(395, 12).Value = 7000
$ws.Cells.Item(395, 13).Value = 7000
$ws.Cells.Item(395, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(395, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(395, 16).Value = 280
$ws.Cells.Item(395, 17).Value = 25
$ws.Cells.Item(395, 18).Value = 'Hortaliza'
$ws.Cells.Item(395, 4).NumberFormat = $ws.Cells.Item(373, 4).NumberFormat

# Dimension will auto-extend to A1:R395 as rows 394-395 now have data